$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily log. It belongs right
# before the (previous) row 127, so insert a fresh row there — this pushes
# every following record down by one (old row 127 -> 128, ... old row
# 240 -> 241) while leaving their data untouched.
$ws.Rows.Item(127).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(127, 1).Value = 10
$ws.Cells.Item(127, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(127, 3).Value = "La Araucanía"
$ws.Cells.Item(127, 4).Value = 44729
$ws.Cells.Item(127, 5).Value = 9
$ws.Cells.Item(127, 6).Value = 100112043
$ws.Cells.Item(127, 7).Value = "Pepino dulce"
$ws.Cells.Item(127, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 35
$ws.Cells.Item(127, 11).Value = 17000
$ws.Cells.Item(127, 12).Value = 17000
$ws.Cells.Item(127, 13).Value = 17000
$ws.Cells.Item(127, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(127, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(127, 16).Value = 944
$ws.Cells.Item(127, 17).Value = 18
$ws.Cells.Item(127, 18).Value = "Hortaliza"
